# Apply the VO terms update described in the commit:
#   - deprecate VO:0021180, replace the reserved-ID placeholder with VO:0021181
#   - keep everything else (including the VO_0005512 - VO_0005560 row) intact
#   - leave the sheet scrolled/selected near the bottom of the list (A18)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell A17 held the next "Reserved IDs" placeholder (VO:0021180); bump it to
# VO:0021181 now that VO:0021180 has been consumed.
$ws.Range("A17").Value = "VO:0021181"

# A2 keeps showing the same "VO_0005512 - VO_0005560" range marker - no text
# change there, just touch it so any shared-string bookkeeping stays tidy.
$ws.Range("A2").Value = "VO_0005512 - VO_0005560"

# Reflect the author's final on-screen selection/scroll position.
$ws.Activate()
$ws.Range("A18").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
